# Update Ultrasonic Sensor Test Results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the two new note strings in column E, next to the top of the table
$ws.Range("E2").Value = "Note* Beyond 1000mm, the distance appears to wrap, indicating a timer overflow."
$ws.Range("E3").Value = " If more than 1m distance required for the project, I'll have to look into timer overflow counting."

# The test data beyond 1000mm (rows 41-52) is no longer valid (timer overflow),
# so clear out the distance values and the computed error formulas but keep
# the existing cell formatting in column A.
$ws.Range("A41:C52").ClearContents() | Out-Null

# Reset the view: scroll back to the top and select the new note cell
$ws.Range("E3").Select() | Out-Null
